# edit.ps1
# Applies the "GITHUB, SERVIDOR REMOTO" addition described by the diff:
#  1. Removes the _GoBack bookmark from its old location (end of the
#     "Estos archivos..." paragraph).
#  2. Inserts two new paragraphs (a heading and a body paragraph) right
#     after the existing "Git log..." paragraph and its two following
#     empty paragraphs, followed by two new empty paragraphs. The
#     _GoBack bookmark is recreated at the end of the new body paragraph,
#     matching Word's normal behavior of tracking the most recent edit.

$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark -------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

# --- Step 2: locate the insertion point -------------------------------
# Find the end of the "Git log: ..." paragraph, then move forward past
# the two empty paragraphs that already follow it, landing right before
# the final (pre-existing) trailing empty paragraph.
$findRange = $d.Content
$found = $findRange.Find.Execute("Git log: sirve para ver el historial de nuestro repositorio.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPoint = $findRange.Duplicate
$insertPoint.Collapse(0)
$insertPoint.Move(4, 3) | Out-Null

# --- Step 3: insert the new paragraphs as OOXML -----------------------
$xmlChunk = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>GITHUB, SERVIDOR REMOTO:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Si deseo pasar los datos que edite desde </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a los archivos que tengo en la carpeta, vamos al control de comandos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>y</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> insertamos el siguiente comando “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pull</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”, automáticamente modificara todos los archivos que tengo en la carpeta de la maquina a lo que modifique en github</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlChunk)
